$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers in I1 and J1 (bold/border style matching the rest of row 1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Fill I2:I31 with 1, and J2:J31 with a copy of the H column values
for ($r = 2; $r -le 31; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hVal
}
